# Add team record (Wins/Losses/Ties) columns to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: AD1="Wins", AE1="Losses", AF1="Ties" - match header style (col A1's style)
$headerSrc = $ws.Cells.Item(1, 1)
$headerSrc.Copy()

$winsHeader = $ws.Cells.Item(1, 30)
$winsHeader.Value = "Wins"
$winsHeader.PasteSpecial(-4122)

$lossesHeader = $ws.Cells.Item(1, 31)
$lossesHeader.Value = "Losses"
$lossesHeader.PasteSpecial(-4122)

$tiesHeader = $ws.Cells.Item(1, 32)
$tiesHeader.Value = "Ties"
$tiesHeader.PasteSpecial(-4122)

# Data rows 2-43: AD = Wins (89), AE = Losses (74), AF = Ties (0)
for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 30).Value = 89
    $ws.Cells.Item($r, 31).Value = 74
    $ws.Cells.Item($r, 32).Value = 0
}

Write-Host "Team record columns added."
